# crisis_sfy_2021.xlsx: fill in the monthly Crisis RPA figures for
# Jul 2020 - Dec 2020 (columns E:J) and the SFY-2021-to-date totals
# (column Q) on the crisis_src sheet. Previously these cells were all 0
# (no data submitted yet); now every month since July 2020 has been run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("crisis_src")

$ws.Range("E3").Value = 262
$ws.Range("F3").Value = 268
$ws.Range("G3").Value = 254
$ws.Range("H3").Value = 264
$ws.Range("I3").Value = 264
$ws.Range("J3").Value = 268
$ws.Range("Q3").Value = 1580

$ws.Range("E4").Value = 30
$ws.Range("F4").Value = 36
$ws.Range("G4").Value = 39
$ws.Range("H4").Value = 45
$ws.Range("I4").Value = 49
$ws.Range("J4").Value = 47
$ws.Range("Q4").Value = 246

$ws.Range("E5").Value = 292
$ws.Range("F5").Value = 304
$ws.Range("G5").Value = 293
$ws.Range("H5").Value = 309
$ws.Range("I5").Value = 313
$ws.Range("J5").Value = 315
$ws.Range("Q5").Value = 1826

$ws.Range("H6").Value = 1
$ws.Range("Q6").Value = 1

$ws.Range("H7").Value = 1
$ws.Range("Q7").Value = 1

$ws.Range("G8").Value = 23
$ws.Range("H8").Value = 53
$ws.Range("I8").Value = 50
$ws.Range("J8").Value = 52
$ws.Range("Q8").Value = 178

$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 3
$ws.Range("J9").Value = 5
$ws.Range("Q9").Value = 13

$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 11
$ws.Range("I10").Value = 9
$ws.Range("J10").Value = 7
$ws.Range("Q10").Value = 28

$ws.Range("G13").Value = 12
$ws.Range("H13").Value = 15
$ws.Range("I13").Value = 15
$ws.Range("J13").Value = 14
$ws.Range("Q13").Value = 56

$ws.Range("G14").Value = 38
$ws.Range("H14").Value = 84
$ws.Range("I14").Value = 77
$ws.Range("J14").Value = 78
$ws.Range("Q14").Value = 277

$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 48
$ws.Range("H17").Value = 16
$ws.Range("I17").Value = 7
$ws.Range("Q17").Value = 93

$ws.Range("E18").Value = 14
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 18
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 13
$ws.Range("Q18").Value = 80

$ws.Range("E19").Value = 12
$ws.Range("F19").Value = 4
$ws.Range("G19").Value = 7
$ws.Range("Q19").Value = 23

$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 9
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 6
$ws.Range("Q20").Value = 38

$ws.Range("E21").Value = 18
$ws.Range("F21").Value = 25
$ws.Range("G21").Value = 26
$ws.Range("H21").Value = 32
$ws.Range("I21").Value = 22
$ws.Range("J21").Value = 31
$ws.Range("Q21").Value = 154

$ws.Range("J22").Value = 6
$ws.Range("Q22").Value = 6

$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 6
$ws.Range("Q23").Value = 7

$ws.Range("J24").Value = 6
$ws.Range("Q24").Value = 6

$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 6
$ws.Range("Q26").Value = 7

$ws.Range("E27").Value = 3
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 7
$ws.Range("Q27").Value = 20

$ws.Range("G28").Value = 11
$ws.Range("H28").Value = 29
$ws.Range("I28").Value = 31
$ws.Range("J28").Value = 27
$ws.Range("Q28").Value = 98

$ws.Range("G29").Value = 7
$ws.Range("H29").Value = 6
$ws.Range("I29").Value = 6
$ws.Range("J29").Value = 11
$ws.Range("Q29").Value = 30
